# Weekly update: insert a new week's Caigua price rows (Primera/Segunda)
# at the top of the data block for "Agrícola del Norte S.A. de Arica",
# pushing the existing rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 47, shifting rows 47:56 down to 49:58.
$ws.Range("A47:A48").EntireRow.Insert()

# New row 47 - "Primera" quality, week of 2021-09-10 (serial 44449)
$ws.Cells.Item(47, 1).Value = 1
$ws.Cells.Item(47, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(47, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(47, 4).Value = 44449
$ws.Cells.Item(47, 5).Value = 15
$ws.Cells.Item(47, 6).Value = 100112036
$ws.Cells.Item(47, 7).Value = "Caigua"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 130
$ws.Cells.Item(47, 11).Value = 6000
$ws.Cells.Item(47, 12).Value = 6500
$ws.Cells.Item(47, 13).Value = 6250
$ws.Cells.Item(47, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(47, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(47, 16).Value = 312
$ws.Cells.Item(47, 17).Value = 20
$ws.Cells.Item(47, 18).Value = "Hortaliza"

# New row 48 - "Segunda" quality, same week
$ws.Cells.Item(48, 1).Value = 1
$ws.Cells.Item(48, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(48, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(48, 4).Value = 44449
$ws.Cells.Item(48, 5).Value = 15
$ws.Cells.Item(48, 6).Value = 100112036
$ws.Cells.Item(48, 7).Value = "Caigua"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Segunda"
$ws.Cells.Item(48, 10).Value = 120
$ws.Cells.Item(48, 11).Value = 4500
$ws.Cells.Item(48, 12).Value = 5000
$ws.Cells.Item(48, 13).Value = 4750
$ws.Cells.Item(48, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(48, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(48, 16).Value = 238
$ws.Cells.Item(48, 17).Value = 20
$ws.Cells.Item(48, 18).Value = "Hortaliza"
